$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.46"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.57%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "29.86"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.13%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.161"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.24%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05746"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.46%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.643"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.45%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.285"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "7.57%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8586"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.12%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8535"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.85%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.43%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07092"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.07%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03146"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "9.88%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09367"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.36%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001530"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.31%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0005967"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-94.22%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006071"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.75%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.530"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.50%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-3.13%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.10%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03320"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.09%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.66%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.493"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "11.07%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.16%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.80%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.37%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004160"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-18.37%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.85%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-25.29%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03748"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.72%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1071"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.13%"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.002199"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.75%"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002949"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-47.52%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009967"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "6.66%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005449"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "7.05%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.08984"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "19.65%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002217"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-19.75%"
